$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet: Version, Date, Contact ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Add a new "Include from FSIII 12" worksheet, cloned from the last
#     "Include from FSIII 11" sheet so it keeps the same layout, column
#     widths and cell styles, then placed right after it. ---
$src = $wb.Worksheets.Item("Include from FSIII 11")
$src.Copy([System.Reflection.Missing]::Value, $src)

$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "Include from FSIII 12"

# Update the UUID value referenced on the new sheet
$newSheet.Range("C2").Value = "aec684bd-c2ea-4ff0-8eb7-6d2cf67fb863"
